$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("marker_template")

$ws.Range("F2:F7").Value = 2

$ws.Activate()
$ws.Range("E6").Select()
